$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This update swaps the data for two pairs of match rows (keeping the
# leading row-index column A untouched), reflecting a re-sync of the
# underlying odds feed where rows 130/131 and 254/256 were reordered.

function Swap-RowData {
    param(
        $ws,
        [int]$row1,
        [int]$row2,
        [string[]]$cols
    )

    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

$cols = @('B','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC','AD')

Swap-RowData $ws 130 131 $cols
Swap-RowData $ws 254 256 $cols
